$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from the existing header cell (H1) onto the new headers
# so they pick up the same bold/border/centered formatting (style index 1)
# instead of Excel minting a brand-new style entry.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for the new columns I (I0) and J (IF), rows 2-20
$values = @(
    @(7, 8),
    @(1, 3),
    @(8, 8),
    @(6, 6),
    @(8, 9),
    @(1, 2),
    @(8, 8),
    @(8, 8),
    @(1, 2),
    @(3, 5),
    @(6, 6),
    @(4, 4),
    @(3, 3),
    @(4, 4),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
